$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# CU-08 "Inscribir alumno" (row 12): mark Estado as "planificado" and Esfuerzo (hrs) as 1
$ws.Range("E12").Value = "planificado"
$ws.Range("F12").Value = 1

# CU-22 "Iniciar sesion" (row 26): mark Estado as "planificado" and Esfuerzo (hrs) as 1
$ws.Range("E26").Value = "planificado"
$ws.Range("F26").Value = 1

# Move the active selection/view to the last edited cell (E26)
$ws.Range("E26").Select()
